$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8337769508361816
$ws.Range("B1").Value = 1.854262113571167
$ws.Range("D1").Value = 1.95109236240387
$ws.Range("E1").Value = 1.034096717834473
